$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- D4: add a hyperlink-styled (but plain-text) URL cell for the 0.1uF capacitor ---
# Copy formatting from an existing "hyperlink style" cell (D2) so the same
# shared cellXf (s="3") is reused instead of creating a brand-new style.
$ws.Range("D2").Copy()
$ws.Range("D4").PasteSpecial(-4122) | Out-Null
$ws.Range("D4").Value = "http://www.mouser.com/ProductDetail/Murata-Electronics/GRM155R61A104KA01D/?qs=sGAEpiMZZMs0AnBnWHyRQEzybnecWqjRhOc0xHRv%252bw8%3d"

# --- New resistor rows (19-25) ---
$ws.Range("A19").Value = "R15,R17,R20,R21,R24"
$ws.Range("B19").Value = " 1.7K RESISTOR"
$ws.Range("D19").Value = "http://www.mouser.com/ProductDetail/Vishay-Thin-Film/PLTT0805Z1721AGT5/?qs=sGAEpiMZZMvdGkrng054t2RPW9MYoEveLabTrIrA%252buo%3d"

$ws.Range("A20").Value = "R16,R18,R19,R22,R23"
$ws.Range("B20").Value = "3.3K RESISTOR"
$ws.Range("D20").Value = "http://www.mouser.com/ProductDetail/Panasonic/ERJ-P6WF3301V/?qs=sGAEpiMZZMvdGkrng054t4TwNrulOmeCGt9o4bxosTAKRbEIXdIPyg%3d%3d"

$ws.Range("A21").Value = "R13"
$ws.Range("B21").Value = "7K RESISTOR"
$ws.Range("D21").Value = "http://www.mouser.com/ProductDetail/Vishay/TNPW08057061BT/?qs=sGAEpiMZZMvdGkrng054twN1Uf5gDWJOW9yH6oo%252bIik%3d"

$ws.Range("A22").Value = "R4,R5,R6,R7,R8,R9,R10,R11,R12,R14"
$ws.Range("B22").Value = "10K RESISTOR"
$ws.Range("D22").Value = "http://www.mouser.com/ProductDetail/Panasonic/ERJ-P6WF1002V/?qs=sGAEpiMZZMvdGkrng054t4TwNrulOmeC3j4fJE09Xf85wPvNPdhA1w%3d%3d"

$ws.Range("A23").Value = "R3,R25,R26,R27"
$ws.Range("B23").Value = "50 RESISTOR"
$ws.Range("D23").Value = "http://www.mouser.com/ProductDetail/Vishay/CRCW080550R0FKTA/?qs=aRXG1QX2Yl9J6LSLM7CpJQ%3d%3d"

$ws.Range("A24").Value = "R1"
$ws.Range("B24").Value = "370 RESISTOR"
$ws.Range("D24").Value = "http://www.mouser.com/ProductDetail/KOA-Speer/RN73H1ETTP3700F10/?qs=sGAEpiMZZMvdGkrng054t%252bj0%252bMDZxyyPf1KQYUrGN6Q%3d"

$ws.Range("A25").Value = "R2"
$ws.Range("B25").Value = "1.5K RESISTOR"
$ws.Range("D25").Value = "http://www.mouser.com/ProductDetail/Panasonic/ERJ-2RKF1501X/?qs=sGAEpiMZZMvdGkrng054t8AJgcdMkx7xOyRawAAbetk%3d"

# --- New connector / misc parts rows (28-30) ---
$ws.Range("B28").Value = "3.5mm Screw Terminal"
$ws.Range("D28").Value = "https://www.sparkfun.com/products/8084"

$ws.Range("B29").Value = "Molex 4 Wire Jumper assembly"
$ws.Range("C29").Value = "Encoder connector @ board"
$ws.Range("D29").Value = "https://www.sparkfun.com/products/9920"

$ws.Range("B30").Value = "Molex 5 Wire jumper assembly"
$ws.Range("C30").Value = "Encoder connector @ encoder"
$ws.Range("D30").Value = "https://www.sparkfun.com/products/9921"

# --- Update the active selection / view to match the edited author's final position ---
$ws.Range("A28").Select() | Out-Null
